$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.162.11"
$ws.Range("E2").Value = "  -2.87%  "

$ws.Range("D3").Value = "1.713.13"
$ws.Range("E3").Value = "  -3.35%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.80%  "

$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4747"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3444"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.14"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07284"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.044"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.69%  "

$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.864"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.02%  "

$ws.Range("D15").Value = "1.711.32"
$ws.Range("E15").Value = "  -3.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.852"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.22%  "

$ws.Range("E18").Value = "  -2.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06359"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.16%  "

$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.604"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "27.197.73"
$ws.Range("E23").Value = "  -2.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.093"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.24%  "

$ws.Range("D28").Value = "1.906.77"
$ws.Range("E28").Value = "  -3.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.085"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.016"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09245"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.589"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.306"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.90%  "

$ws.Range("E35").Value = "  -4.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05901"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2014"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.749"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.61%  "

$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.411"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.96%  "

$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5931"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.11%  "

$ws.Range("E43").Value = "  -5.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.476"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.570"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5630"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.48%  "

$ws.Range("E49").Value = "  -5.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06632"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.75%  "

$ws.Range("E51").Value = "  -4.77%  "

